# Add a new "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell formatting (bold font, thin border, centered) from
# the existing "sum" header (G1) onto the new H1 header cell, then set
# its text. Using PasteSpecial(formats) reuses the existing header style
# instead of creating a brand new, duplicate style entry.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"

# Fill in the "Save" flag values for the data rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
